$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.156.70"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.832.56"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'241.52"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'0.6648"
$ws.Range("E6").Value = "  -2.33%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.07421"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "'0.2939"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").Value = "'22.79"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").Value = "'0.07749"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "1.838.85"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'4.986"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Value = "'0.6686"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "'82.85"
$ws.Range("E15").Value = "  -4.58%  "
$ws.Range("D16").Value = "'6.094"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "'0.000008347"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "29.141.77"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "2.080.19"
$ws.Range("D20").Value = "'229.09"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'12.47"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D23").Value = "'7.163"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("D25").Value = "'159.50"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").Value = "'0.1404"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("D27").Value = "'8.606"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").Value = "'18.02"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "'1.510"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "'4.112"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").Value = "'4.037"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "'0.05318"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "'1.868"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'0.7487"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").Value = "'1.137"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "'2.641"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").Value = "1.272.92"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").Value = "'0.01797"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").Value = "'2.732"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "'0.9286"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "'0.08506"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("D43").Value = "'5.901"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'101.93"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").Value = "1.986.23"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'0.5147"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").Value = "'1.763"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").Value = "'63.04"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  -0.87%  "
